# Update "想去人数" (F column) figures that changed between scrapes.
# Sheet "展览" (sheet1) and sheet "全部类型" (sheet4) both list the same
# events (sheet4 has one extra row), so both need the matching updates.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetAllTypes   = $wb.Worksheets.Item("全部类型")

# Row -> new F value for "展览"
$updatesExhibition = @{
    5  = 1290
    6  = 17973
    10 = 6779
    11 = 682
    16 = 12
    19 = 210
    26 = 977
    27 = 112
    28 = 5160
    31 = 12
    33 = 12001
    36 = 201
    37 = 269
    39 = 298
    40 = 90
}

foreach ($row in $updatesExhibition.Keys) {
    $sheetExhibition.Cells.Item($row, 6).Value = $updatesExhibition[$row]
}

# Row -> new F value for "全部类型"
$updatesAllTypes = @{
    5  = 1290
    6  = 17973
    10 = 6779
    11 = 682
    16 = 12
    19 = 210
    26 = 977
    27 = 112
    28 = 5160
    33 = 12
    35 = 12001
    38 = 201
    39 = 269
    41 = 298
    42 = 90
}

foreach ($row in $updatesAllTypes.Keys) {
    $sheetAllTypes.Cells.Item($row, 6).Value = $updatesAllTypes[$row]
}
